# "upgrade left table until javakheti"
# Adds a 2023 column (K) to the Sighnaghi average-monthly-remuneration table,
# mirroring the formatting of the existing 2022 column (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles/number formats/borders) from column J's data rows
# down into column K before writing the new values, so the new cells inherit
# the same look (right-aligned custom number format, borders, etc.) as the
# rest of the table.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)

# New 2023 year header
$ws.Range("K3").Value = 2023

# New 2023 data points
$ws.Range("K4").Value = 827.2
$ws.Range("K5").Value = 471.4
$ws.Range("K6").Value = 1221
